# Apply targeted updates to the generated Weekly Units Completed report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Report Generated On" timestamp.
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:02 AM"

# 2. Update Total Billed Amount summary figure.
$ws.Range("C8").Value = 84.24

# 3. Clear the Scope ID # value (was "#NO MATCH", now blank).
$ws.Range("G10").Value = ""

# 4. Update the line item pricing total.
$ws.Range("H16").Value = 84.24

# 5. Update the grand TOTAL pricing figure.
$ws.Range("H17").Value = 84.24
